# Apply cell value updates for the cryptos worksheet refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.784.32"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.346.37"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'239.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").Value = "'73.21"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Value = "'58.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("D12").Value = "'32.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").Value = "'7.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "2.696.93"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "'16.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").Value = "'0.903"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "2.348.32"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "43.721.79"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "'0.0000102"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "'6.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "'77.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'256.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  +22.48%  "
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'2.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").Value = "'22.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'177.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'5.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'5.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("D37").Value = "'3.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("D38").Value = "'2.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("D39").Value = "'6.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.20%  "
$ws.Range("D40").Value = "'0.0280"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("D41").Value = "'69.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +32.01%  "
$ws.Range("E42").Value = "  +11.70%  "
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").Value = "'19.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("E45").Value = "  +7.10%  "
$ws.Range("D46").Value = "'4.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.19%  "
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").Value = "'99.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "'1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.95%  "
